# Generate Report for Handback
# For both the "zh-cn" and "de-de" localization sheets, mark rows 2 and 3
# as handed back: update the Status text, fill in the "Latest Target File"
# and "Latest Handback File" hyperlink columns (F, G) with the same file
# names used for handoff, and stamp the "Latest Handback DateTime" column
# (H) with the real handback timestamp.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# ---- zh-cn sheet ----

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/855e21738df730ae835f3bc3473e218d1ba2d427/e2e/0d6b0f35-16cc-4cb8-bd66-4e305df6314e.md", "", "", "0d6b0f35-16cc-4cb8-bd66-4e305df6314e.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8d6acc2002b8c5dd28805d38cd4320620172bf2f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0d6b0f35-16cc-4cb8-bd66-4e305df6314e.142b8a8ba9b052f42730fd821938c0c609886cb5.zh-cn.xlf", "", "", "0d6b0f35-16cc-4cb8-bd66-4e305df6314e.142b8a8ba9b052f42730fd821938c0c609886cb5.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/855e21738df730ae835f3bc3473e218d1ba2d427/e2e/97814758-b34a-46ff-8abf-bdbd317417c5.md", "", "", "97814758-b34a-46ff-8abf-bdbd317417c5.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8d6acc2002b8c5dd28805d38cd4320620172bf2f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.zh-cn.xlf", "", "", "97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.zh-cn.xlf")

$zhcn.Range("H2").Value = "2016-03-18 05:46:38"
$zhcn.Range("H3").Value = "2016-03-18 05:46:38"

# ---- de-de sheet ----

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/855e21738df730ae835f3bc3473e218d1ba2d427/e2e/0d6b0f35-16cc-4cb8-bd66-4e305df6314e.md", "", "", "0d6b0f35-16cc-4cb8-bd66-4e305df6314e.md")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8bf56e1405336e77427f9587d4be0b023c428ba2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0d6b0f35-16cc-4cb8-bd66-4e305df6314e.142b8a8ba9b052f42730fd821938c0c609886cb5.de-de.xlf", "", "", "0d6b0f35-16cc-4cb8-bd66-4e305df6314e.142b8a8ba9b052f42730fd821938c0c609886cb5.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/855e21738df730ae835f3bc3473e218d1ba2d427/e2e/97814758-b34a-46ff-8abf-bdbd317417c5.md", "", "", "97814758-b34a-46ff-8abf-bdbd317417c5.md")
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8bf56e1405336e77427f9587d4be0b023c428ba2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.de-de.xlf", "", "", "97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.de-de.xlf")

$dede.Range("H2").Value = "2016-03-18 05:46:43"
$dede.Range("H3").Value = "2016-03-18 05:46:43"
